# Edit described by the commit:
#   1. Slide 5's table switches from the custom "Table_0" style to the
#      built-in table style {6EECB4D3-6956-4014-A919-41729CB0C93B}.
#   2. The presentation's theme colour scheme changes from the "Integral"
#      (Red Violet) palette to the default "Office Theme" palette - i.e.
#      the deck's Design/Theme was switched in the UI.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
# Find the (single) table in the deck - it lives on slide 5 - and swap it
# from the custom "Table_0" style to the built-in style the author picked
# from the Table Design gallery.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{6EECB4D3-6956-4014-A919-41729CB0C93B}")
        }
    }
}

# --- 2. Theme colours ------------------------------------------------------
# New ("Office Theme") RGB values, converted to the BGR-packed long that the
# ColorScheme/RGB COM property expects, in the fixed
# dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink order.
$officeThemeColors = @(
    0,            # dk1     000000
    16777215,     # lt1     FFFFFF
    6968388,      # dk2     44546A
    15132391,     # lt2     E7E6E6
    13998939,     # accent1 5B9BD5
    3243501,      # accent2 ED7D31
    10855845,     # accent3 A5A5A5
    49407,        # accent4 FFC000
    12874308,     # accent5 4472C4
    4697456,      # accent6 70AD47
    12673797,     # hlink   0563C1
    7491477       # folHlink 954F72
)

$colorScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
